# Convert the color-wise fraction values (B2:J17) into percentage values
# by multiplying each numeric cell by 100 and rounding to 2 decimal places,
# then set the explicit column widths that Excel recorded for this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 17; $r++) {
    for ($c = 2; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value2
        if ($old -ne $null) {
            $new = [Math]::Round($old * 100, 2)
            $cell.Value = $new
        }
    }
}

# Explicit column widths (in Excel "characters" units) recorded on the sheet.
$ws.Columns.Item(1).ColumnWidth = 17.09
$ws.Columns.Item(2).ColumnWidth = 9.92
$ws.Columns.Item(3).ColumnWidth = 7.42
$ws.Columns.Item(4).ColumnWidth = 8.75
$ws.Columns.Item(5).ColumnWidth = 8.75
$ws.Columns.Item(6).ColumnWidth = 7.42
$ws.Columns.Item(7).ColumnWidth = 5.09
$ws.Columns.Item(8).ColumnWidth = 14.75
$ws.Columns.Item(9).ColumnWidth = 6.25
$ws.Columns.Item(10).ColumnWidth = 59.09
